{"js": "// Replace each three-digit-number division expression in the document\n// with its new value, per the commit's regenerated numbers.\nconst replacements = [\n  [\"255\u00f77=\", \"442\u00f73=\"],\n  [\"287\u00f78=\", \"826\u00f78=\"],\n  [\"487\u00f77=\", \"953\u00f73=\"],\n  [\"507\u00f75=\", \"817\u00f73=\"],\n  [\"445\u00f73=\", \"518\u00f77=\"],\n  [\"475\u00f75=\", \"960\u00f77=\"],\n  [\"167\u00f76=\", \"390\u00f73=\"],\n  [\"427\u00f72=\", \"961\u00f72=\"],\n  [\"978\u00f72=\", \"780\u00f78=\"],\n  [\"517\u00f77=\", \"988\u00f79=\"],\n  [\"283\u00f78=\", \"469\u00f78=\"],\n  [\"158\u00f75=\", \"394\u00f73=\"],\n  [\"322\u00f73=\", \"212\u00f79=\"],\n  [\"763\u00f72=\", \"983\u00f79=\"],\n  [\"768\u00f77=\", \"942\u00f75=\"],\n  [\"976\u00f74=\", \"473\u00f78=\"],\n  [\"914\u00f76=\", \"167\u00f74=\"],\n  [\"930\u00f79=\", \"185\u00f79=\"],\n  [\"198\u00f73=\", \"805\u00f76=\"],\n  [\"702\u00f75=\", \"503\u00f75=\"],\n  [\"430\u00f78=\", \"623\u00f74=\"],\n  [\"702\u00f78=\", \"197\u00f79=\"],\n  [\"418\u00f77=\", \"156\u00f74=\"],\n  [\"869\u00f73=\", \"351\u00f79=\"],\n  [\"663\u00f75=\", \"343\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-number division expression in the document\n# with its new value, per the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"255\u00f77=\", \"442\u00f73=\"),\n    @(\"287\u00f78=\", \"826\u00f78=\"),\n    @(\"487\u00f77=\", \"953\u00f73=\"),\n    @(\"507\u00f75=\", \"817\u00f73=\"),\n    @(\"445\u00f73=\", \"518\u00f77=\"),\n    @(\"475\u00f75=\", \"960\u00f77=\"),\n    @(\"167\u00f76=\", \"390\u00f73=\"),\n    @(\"427\u00f72=\", \"961\u00f72=\"),\n    @(\"978\u00f72=\", \"780\u00f78=\"),\n    @(\"517\u00f77=\", \"988\u00f79=\"),\n    @(\"283\u00f78=\", \"469\u00f78=\"),\n    @(\"158\u00f75=\", \"394\u00f73=\"),\n    @(\"322\u00f73=\", \"212\u00f79=\"),\n    @(\"763\u00f72=\", \"983\u00f79=\"),\n    @(\"768\u00f77=\", \"942\u00f75=\"),\n    @(\"976\u00f74=\", \"473\u00f78=\"),\n    @(\"914\u00f76=\", \"167\u00f74=\"),\n    @(\"930\u00f79=\", \"185\u00f79=\"),\n    @(\"198\u00f73=\", \"805\u00f76=\"),\n    @(\"702\u00f75=\", \"503\u00f75=\"),\n    @(\"430\u00f78=\", \"623\u00f74=\"),\n    @(\"702\u00f78=\", \"197\u00f79=\"),\n    @(\"418\u00f77=\", \"156\u00f74=\"),\n    @(\"869\u00f73=\", \"351\u00f79=\"),\n    @(\"663\u00f75=\", \"343\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
